$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bump the embedded iframe size from 640x360 to 1280x720 for every row.
# Column D holds the shared trailing-HTML string used by the CONCATENATE
# formula in column B, so updating it updates every row at once.
$newSuffix = '" width=1280" height="720" frameborder="0" webkitallowfullscreen mozallowfullscreen allowfullscreen></iframe>'
for ($r = 2; $r -le 19; $r++) {
    $ws.Range("D$r").Value = $newSuffix
}

# Fill in the previously-missing Vimeo source URLs in column F so that the
# CONCATENATE/LEFT/RIGHT formulas in columns B and E stop erroring out and
# resolve to real values ("16 View Videos Completed").
$ws.Range("F4").Value  = "https://vimeo.com/175554706/fe8a1a8ef2"
$ws.Range("F5").Value  = "https://vimeo.com/175554705/400e226eb3"
$ws.Range("F6").Value  = "https://vimeo.com/175849908/00ec6c6867"
$ws.Range("F7").Value  = "https://vimeo.com/175850249/cda867f897"
$ws.Range("F9").Value  = "https://vimeo.com/175554709/b03e806831"
$ws.Range("F10").Value = "https://vimeo.com/175554708/4647d2b2db"
$ws.Range("F11").Value = "https://vimeo.com/175554710/b1d449bed5"
$ws.Range("F18").Value = "https://vimeo.com/176225772/6441d88ef9"
$ws.Range("F19").Value = "https://vimeo.com/176226121/d76c0f4be8"

# Recalculate so the dependent formulas (B2:B19, E2:E19) pick up the new
# source data instead of keeping stale cached #VALUE! results.
$excel.Calculate()

# Cosmetic: update the saved selection to match the author's final state.
$ws.Range("B19").Select()
